$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = New-Object 'object[,]' 1,10
$row[0,0] = 'فصل سوم منتهی به 1399/09'
$row[0,1] = 'فصل چهارم منتهی به 1399/12'
$row[0,2] = 'فصل اول منتهی به 1400/03'
$row[0,3] = 'فصل دوم منتهی به 1400/06'
$row[0,4] = 'فصل سوم منتهی به 1400/09'
$row[0,5] = 'فصل چهارم منتهی به 1400/12'
$row[0,6] = 'فصل اول منتهی به 1401/03'
$row[0,7] = 'فصل دوم منتهی به 1401/06'
$row[0,8] = 'فصل سوم منتهی به 1401/09'
$row[0,9] = 'فصل چهارم منتهی به 1401/12'
$ws.Range("D8:M8").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = '1400-10-30 (2)'
$row[0,1] = '1401-03-04 (8)'
$row[0,2] = '1401-04-30 (2)'
$row[0,3] = '1401-08-30 (4)'
$row[0,4] = '1401-10-28 (2)'
$row[0,5] = '1402-02-28 (7)'
$row[0,6] = '1401-04-30'
$row[0,7] = '1401-08-30 (2)'
$row[0,8] = '1401-10-28'
$row[0,9] = '1402-02-28'
$ws.Range("D9:M9").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -487686
$row[0,1] = 529080
$row[0,2] = 347296
$row[0,3] = -813518
$row[0,4] = 220415
$row[0,5] = 2227690
$row[0,6] = 1209295
$row[0,7] = 1350666
$row[0,8] = 125319
$row[0,9] = 1014673
$ws.Range("D12:M12").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 5755
$row[0,2] = -20000
$row[0,3] = -19646
$row[0,4] = -17000
$row[0,5] = -11965
$row[0,6] = -21462
$row[0,7] = -366214
$row[0,8] = -235968
$row[0,9] = -220066
$ws.Range("D13:M13").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -487686
$row[0,1] = 534835
$row[0,2] = 327296
$row[0,3] = -833164
$row[0,4] = 203415
$row[0,5] = 2215725
$row[0,6] = 1187833
$row[0,7] = 984452
$row[0,8] = -110649
$row[0,9] = 794607
$ws.Range("D14:M14").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D16:M16").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -97136
$row[0,1] = -174078
$row[0,2] = -81195
$row[0,3] = -414386
$row[0,4] = -269657
$row[0,5] = -381832
$row[0,6] = -450729
$row[0,7] = -343107
$row[0,8] = -268415
$row[0,9] = -679120
$ws.Range("D17:M17").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D18:M18").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D19:M19").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = -2559
$row[0,4] = 0
$row[0,5] = -296
$row[0,6] = 0
$row[0,7] = -1114
$row[0,8] = 0
$row[0,9] = -1238
$ws.Range("D20:M20").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D21:M21").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D22:M22").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = -5668
$row[0,2] = 0
$row[0,3] = 179241
$row[0,4] = 0
$row[0,5] = 340679
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 37679
$row[0,9] = 131390
$ws.Range("D23:M23").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 116
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D24:M24").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D25:M25").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = -100000
$ws.Range("D26:M26").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D27:M27").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D28:M28").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D29:M29").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D30:M30").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 19
$row[0,1] = 12
$row[0,2] = 0
$row[0,3] = 121
$row[0,4] = -121
$row[0,5] = 2804
$row[0,6] = 0
$row[0,7] = 64
$row[0,8] = 979
$row[0,9] = 660
$ws.Range("D31:M31").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -97117
$row[0,1] = -179618
$row[0,2] = -81195
$row[0,3] = -237583
$row[0,4] = -269778
$row[0,5] = -38645
$row[0,6] = -450729
$row[0,7] = -344157
$row[0,8] = -229757
$row[0,9] = -648308
$ws.Range("D32:M32").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -584803
$row[0,1] = 355217
$row[0,2] = 246101
$row[0,3] = -1070747
$row[0,4] = -66363
$row[0,5] = 2177080
$row[0,6] = 737104
$row[0,7] = 640295
$row[0,8] = -340406
$row[0,9] = 146299
$ws.Range("D33:M33").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 1842
$row[0,2] = 0
$row[0,3] = 88361
$row[0,4] = 0
$row[0,5] = -69038
$row[0,6] = 65550
$row[0,7] = -65550
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D35:M35").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = '-'
$row[0,1] = '-'
$row[0,2] = '-'
$row[0,3] = '-'
$row[0,4] = '-'
$row[0,5] = 0
$row[0,6] = '-'
$row[0,7] = '-'
$row[0,8] = '-'
$row[0,9] = 0
$ws.Range("D36:M36").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D37:M37").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = -25000
$row[0,4] = 25000
$row[0,5] = -46000
$row[0,6] = -21500
$row[0,7] = -978
$row[0,8] = 0
$row[0,9] = -42260
$ws.Range("D38:M38").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 5349957
$row[0,1] = 3633014
$row[0,2] = 3548504
$row[0,3] = 6293746
$row[0,4] = 5751600
$row[0,5] = 4407304
$row[0,6] = 6071376
$row[0,7] = 5517524
$row[0,8] = 5219259
$row[0,9] = 8919748
$ws.Range("D39:M39").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -4331939
$row[0,1] = -3593744
$row[0,2] = -3347734
$row[0,3] = -4235968
$row[0,4] = -5333430
$row[0,5] = -4463911
$row[0,6] = -5476960
$row[0,7] = -6477339
$row[0,8] = -4267019
$row[0,9] = -6430458
$ws.Range("D40:M40").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -220047
$row[0,1] = -83216
$row[0,2] = -247374
$row[0,3] = -214567
$row[0,4] = -91257
$row[0,5] = -847677
$row[0,6] = -362632
$row[0,7] = -618620
$row[0,8] = -445041
$row[0,9] = -478336
$ws.Range("D41:M41").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D42:M42").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D43:M43").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D44:M44").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D45:M45").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D46:M46").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D47:M47").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D48:M48").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D49:M49").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -231554
$row[0,1] = -327835
$row[0,2] = -252155
$row[0,3] = -392704
$row[0,4] = -314347
$row[0,5] = -836653
$row[0,6] = -388621
$row[0,7] = -2394
$row[0,8] = -253825
$row[0,9] = -1110865
$ws.Range("D50:M50").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 566417
$row[0,1] = -369939
$row[0,2] = -298759
$row[0,3] = 1513868
$row[0,4] = 37566
$row[0,5] = -1855975
$row[0,6] = -112787
$row[0,7] = -1647357
$row[0,8] = 253374
$row[0,9] = 857829
$ws.Range("D51:M51").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = -18386
$row[0,1] = -14722
$row[0,2] = -52658
$row[0,3] = 443121
$row[0,4] = -28797
$row[0,5] = 321105
$row[0,6] = 624317
$row[0,7] = -1007062
$row[0,8] = -87032
$row[0,9] = 1004128
$ws.Range("D52:M52").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 353606
$row[0,1] = 335220
$row[0,2] = 320498
$row[0,3] = 267840
$row[0,4] = 710961
$row[0,5] = 682164
$row[0,6] = 1003269
$row[0,7] = 1627586
$row[0,8] = 620524
$row[0,9] = 533492
$ws.Range("D53:M53").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 0
$row[0,2] = 0
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D54:M54").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 335220
$row[0,1] = 320498
$row[0,2] = 267840
$row[0,3] = 710961
$row[0,4] = 682164
$row[0,5] = 1003269
$row[0,6] = 1627586
$row[0,7] = 620524
$row[0,8] = 533492
$row[0,9] = 1537620
$ws.Range("D55:M55").Value = $row

$row = New-Object 'object[,]' 1,10
$row[0,0] = 0
$row[0,1] = 878158
$row[0,2] = 0
$row[0,3] = 78940
$row[0,4] = -78940
$row[0,5] = 1480677
$row[0,6] = 0
$row[0,7] = 0
$row[0,8] = 0
$row[0,9] = 0
$ws.Range("D56:M56").Value = $row
